$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Points = 2)
$ws.Range("B3").Value = -3.33225
$ws.Range("C3").Value = -1.7599
$ws.Range("D3").Value = 2.13937
$ws.Range("E3").Value = "Fail"

# Row 4 (Points = 3)
$ws.Range("B4").Value = 0.11167
$ws.Range("C4").Value = -0.13056
$ws.Range("D4").Value = -0.466
$ws.Range("E4").Value = "Pass"

# Row 5 (Points = 4)
$ws.Range("B5").Value = 1.84916
$ws.Range("C5").Value = 2.55912
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = "Fail"

# Row 6 (Points = 5)
$ws.Range("B6").Value = 0.8031
$ws.Range("C6").Value = 0.89117
$ws.Range("D6").Value = -0.04764

# Row 7 (Points = 6)
$ws.Range("B7").Value = -0.82789
$ws.Range("C7").Value = 15.85663
$ws.Range("D7").Value = -1.37523
$ws.Range("E7").Value = "Fail"

# Row 8 (Points = 7)
$ws.Range("B8").Value = 0.34998
$ws.Range("C8").Value = 0.28921
$ws.Range("D8").Value = -1.36607

# Row 9 (Points = 8)
$ws.Range("B9").Value = -1.5188
$ws.Range("C9").Value = -1.61244
$ws.Range("D9").Value = -1.66409

# Row 10 (Points = 9)
$ws.Range("B10").Value = -1.39746
$ws.Range("C10").Value = -1.13727
$ws.Range("D10").Value = -1.18109

# Row 11 (Points = 10)
$ws.Range("B11").Value = -0.05815
$ws.Range("C11").Value = 0.3109
$ws.Range("D11").Value = -0.75213
